$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-7 from 2023-09-16 (45185) to 2023-10-05 (45204)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45204
}
